$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.870.94'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '1.833.05'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.10'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("E6").Value = '  -2.08%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07710'
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3049'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.38'
$ws.Range("E10").Value = '  -4.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07806'
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '1.827.36'
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.092'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.54'
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6820'
$ws.Range("E15").Value = '  -2.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.414'
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008308'
$ws.Range("E17").Value = '  -3.03%  '
$ws.Range("D18").Value = '28.875.06'
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.43'
$ws.Range("E19").Value = '  -3.14%  '
$ws.Range("D20").Value = '2.075.70'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.452'
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1477'
$ws.Range("E25").Value = '  -3.99%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.798'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.05'
$ws.Range("E27").Value = '  -1.83%  '
$ws.Range("E28").Value = '  -2.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.544'
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.224'
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.152'
$ws.Range("E31").Value = '  -2.58%  '
$ws.Range("E32").Value = '  -0.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05094'
$ws.Range("E33").Value = '  -3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7771'
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.851'
$ws.Range("E35").Value = '  -2.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.141'
$ws.Range("E36").Value = '  -3.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.690'
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01852'
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").Value = '1.225.96'
$ws.Range("E39").Value = '  -3.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.693'
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9546'
$ws.Range("E41").Value = '  +6.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.20'
$ws.Range("E42").Value = '  -0.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.878'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.623'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("E46").Value = '  -3.69%  '
$ws.Range("D47").Value = '1.976.84'
$ws.Range("E47").Value = '  -2.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5156'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.12'
$ws.Range("E49").Value = '  -9.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.751'
$ws.Range("E50").Value = '  -2.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.923'
$ws.Range("E51").Value = '  -1.79%  '
